# Insert a new data row at row 158 (pushing existing row 158 and all
# subsequent rows down by one), then populate the new row with the
# values describing the new weekly record. All of the "descriptive"
# (non price/date) columns are copied from the row directly below
# (which used to be row 158 before the insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 158 downward (and everything after it) by inserting a new row.
$ws.Rows.Item(158).Insert()

# The row that used to be 158 is now row 159; copy its "template" values
# (the columns that do not change for the new record) into the new row 158.
$srcRow = 159
$dstRow = 158

$ws.Cells.Item($dstRow, 1).Value = $ws.Cells.Item($srcRow, 1).Value2   # A Mercado ID
$ws.Cells.Item($dstRow, 2).Value = $ws.Cells.Item($srcRow, 2).Value2   # B Mercado
$ws.Cells.Item($dstRow, 3).Value = $ws.Cells.Item($srcRow, 3).Value2   # C Region

$ws.Cells.Item($dstRow, 4).Value = 44729                               # D Fecha

$ws.Cells.Item($dstRow, 5).Value = $ws.Cells.Item($srcRow, 5).Value2   # E Codreg
$ws.Cells.Item($dstRow, 6).Value = $ws.Cells.Item($srcRow, 6).Value2   # F Categoria ID
$ws.Cells.Item($dstRow, 7).Value = $ws.Cells.Item($srcRow, 7).Value2   # G Categoria
$ws.Cells.Item($dstRow, 8).Value = $ws.Cells.Item($srcRow, 8).Value2   # H Variedad
$ws.Cells.Item($dstRow, 9).Value = $ws.Cells.Item($srcRow, 9).Value2   # I Calidad

$ws.Cells.Item($dstRow, 10).Value = 400                                # J Volumen
$ws.Cells.Item($dstRow, 11).Value = 22000                              # K Precio minimo
$ws.Cells.Item($dstRow, 12).Value = 23000                              # L Precio maximo
$ws.Cells.Item($dstRow, 13).Value = 22500                              # M Precio promedio ponderado

$ws.Cells.Item($dstRow, 14).Value = $ws.Cells.Item($srcRow, 14).Value2 # N Unidad de comercializacion
$ws.Cells.Item($dstRow, 15).Value = $ws.Cells.Item($srcRow, 15).Value2 # O Origen

$ws.Cells.Item($dstRow, 16).Value = 375                                # P Precio $/Kg

$ws.Cells.Item($dstRow, 17).Value = $ws.Cells.Item($srcRow, 17).Value2 # Q Kg o Unidades
$ws.Cells.Item($dstRow, 18).Value = $ws.Cells.Item($srcRow, 18).Value2 # R Clasificacion

# Match the date number format used by the other rows in column D.
$ws.Cells.Item($dstRow, 4).NumberFormat = $ws.Cells.Item($srcRow, 4).NumberFormat
